$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("master-reg_center_user")

# Append three new rows of test data (16th May refresh) after the last
# existing row (row 33), mirroring the existing rows' layout/types.
$newRows = @(
    @(10005, 110033),
    @(10005, 110034),
    @(10005, 110035)
)

$startRow = 34
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $regcntrId = $newRows[$i][0]
    $usrId = $newRows[$i][1]

    $ws.Cells.Item($r, 1).Value = $regcntrId
    $ws.Cells.Item($r, 2).Value = $usrId
    $ws.Cells.Item($r, 3).Value = "eng"
    $ws.Cells.Item($r, 4).Value = $true
    $ws.Cells.Item($r, 5).Value = "superadmin"
    $ws.Cells.Item($r, 6).Value = "now()"
}

# Select the row below the newly-added data, matching the post-edit
# workbook state (full-row selection of row 37 downward), and scroll the
# view back to the top-left.
$ws.Range("A37:XFD1048576").Select()
$excel.ActiveWindow.ScrollRow = 1
